$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data refresh: update Fecha (D), Volumen (J), Precio minimo (K),
# Precio maximo (L), Precio promedio ponderado (M) and Precio $/Kg (P)
# for each data row (2-35) with the latest week's figures.

$rows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35)
$dVals = @(44165, 44193, 44627, 44537, 44260, 44804, 44225, 44223, 44166, 44179, 44340, 44242, 44536, 44221, 44389, 44935, 44669, 44390, 44574, 44965, 44243, 45092, 44224, 44845, 44992, 44557, 44222, 44967, 44292, 44187, 44937, 44756, 44559, 44291)
$jVals = @(68, 70, 78, 88, 60, 85, 56, 80, 45, 78, 54, 95, 125, 50, 81, 78, 92, 50, 50, 87, 45, 90, 67, 80, 45, 104, 45, 110, 40, 65, 68, 104, 68, 45)
$kVals = @(3000, 3000, 3500, 2000, 3500, 3000, 3000, 2500, 2500, 3000, 3000, 2500, 2200, 2500, 2800, 3000, 2500, 3000, 3000, 3000, 3000, 3000, 3000, 2500, 4000, 2000, 3000, 3000, 3000, 3000, 3500, 2800, 2000, 3000)
$lVals = @(3000, 3000, 3500, 2200, 3500, 3000, 3000, 3000, 2500, 3000, 3000, 3000, 2200, 2500, 3000, 3000, 3000, 3000, 3000, 3000, 3000, 3500, 3000, 2500, 4000, 2500, 3000, 3300, 3000, 3000, 3500, 3000, 2000, 3000)
$mVals = @(3000, 3000, 3500, 2091, 3500, 3000, 3000, 2781, 2500, 3000, 3000, 2737, 2200, 2500, 2889, 3000, 2755, 3000, 3000, 3000, 3000, 3278, 3000, 2500, 4000, 2260, 3000, 3136, 3000, 3000, 3500, 2904, 2000, 3000)
$pVals = @(1000, 1000, 1167, 697, 1167, 1000, 1000, 927, 833, 1000, 1000, 912, 733, 833, 963, 1000, 918, 1000, 1000, 1000, 1000, 1093, 1000, 833, 1333, 753, 1000, 1045, 1000, 1000, 1167, 968, 667, 1000)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 4).Value2  = $dVals[$i]   # D: Fecha
    $ws.Cells.Item($r, 10).Value2 = $jVals[$i]   # J: Volumen
    $ws.Cells.Item($r, 11).Value2 = $kVals[$i]   # K: Precio minimo
    $ws.Cells.Item($r, 12).Value2 = $lVals[$i]   # L: Precio maximo
    $ws.Cells.Item($r, 13).Value2 = $mVals[$i]   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value2 = $pVals[$i]   # P: Precio $/Kg
}
